$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

$ws.Cells.Item(2, 2).Value = 524.25
$ws.Cells.Item(2, 3).Value = 520.85
$ws.Cells.Item(3, 2).Value = 8783
$ws.Cells.Item(3, 3).Value = 8690.700000000001
$ws.Cells.Item(4, 2).Value = 3002
$ws.Cells.Item(4, 3).Value = 3013.75
$ws.Cells.Item(5, 2).Value = 505.75
$ws.Cells.Item(5, 3).Value = 504.25
$ws.Cells.Item(6, 2).Value = 223.28
$ws.Cells.Item(6, 3).Value = 219.56
$ws.Cells.Item(7, 2).Value = 1819.9
$ws.Cells.Item(7, 3).Value = 1815.8
$ws.Cells.Item(8, 2).Value = 6899.55
$ws.Cells.Item(8, 3).Value = 6899.5
$ws.Cells.Item(9, 2).Value = 191.96
$ws.Cells.Item(9, 3).Value = 190.19
$ws.Cells.Item(10, 2).Value = 253.85
$ws.Cells.Item(10, 3).Value = 254.2
$ws.Cells.Item(11, 2).Value = 247.97
$ws.Cells.Item(11, 3).Value = 241.96
$ws.Cells.Item(12, 2).Value = 52310.4
$ws.Cells.Item(12, 3).Value = 51491.7
$ws.Cells.Item(13, 2).Value = 15381.8
$ws.Cells.Item(13, 3).Value = 15130.85
$ws.Cells.Item(14, 2).Value = 875.15
$ws.Cells.Item(14, 3).Value = 861
$ws.Cells.Item(15, 2).Value = 4765.65
$ws.Cells.Item(15, 3).Value = 4623.85
$ws.Cells.Item(16, 2).Value = 3815.4
$ws.Cells.Item(16, 3).Value = 3789
$ws.Cells.Item(17, 2).Value = 195.3
$ws.Cells.Item(17, 3).Value = 193.61
$ws.Cells.Item(18, 2).Value = 1857.7
$ws.Cells.Item(18, 3).Value = 1867.8
$ws.Cells.Item(19, 2).Value = 753.5
$ws.Cells.Item(19, 3).Value = 734.8
$ws.Cells.Item(20, 2).Value = 451.7
$ws.Cells.Item(20, 3).Value = 504.55
$ws.Cells.Item(21, 2).Value = 1347.25
$ws.Cells.Item(21, 3).Value = 1347.35
$ws.Cells.Item(22, 2).Value = 964.5
$ws.Cells.Item(22, 3).Value = 930
$ws.Cells.Item(23, 2).Value = 618.6
$ws.Cells.Item(23, 3).Value = 611.8
$ws.Cells.Item(24, 2).Value = 2964.25
$ws.Cells.Item(24, 3).Value = 2964.6
$ws.Cells.Item(25, 2).Value = 291
$ws.Cells.Item(25, 3).Value = 282.15
$ws.Cells.Item(26, 2).Value = 24949.15
$ws.Cells.Item(26, 3).Value = 24843.8
$ws.Cells.Item(27, 2).Value = 424.95
$ws.Cells.Item(27, 3).Value = 417.75
$ws.Cells.Item(28, 2).Value = 283.3
$ws.Cells.Item(28, 3).Value = 281.6
$ws.Cells.Item(29, 2).Value = 546.25
$ws.Cells.Item(29, 3).Value = 542.35
$ws.Cells.Item(30, 2).Value = 820.4
$ws.Cells.Item(30, 3).Value = 811.05
$ws.Cells.Item(31, 2).Value = 759.95
$ws.Cells.Item(31, 3).Value = 765.1
$ws.Cells.Item(32, 2).Value = 910.15
$ws.Cells.Item(32, 3).Value = 891.6
$ws.Cells.Item(33, 2).Value = 453.55
$ws.Cells.Item(33, 3).Value = 450.2
$ws.Cells.Item(34, 2).Value = 155.39
$ws.Cells.Item(34, 3).Value = 152.4
$ws.Cells.Item(35, 2).Value = 480.85
$ws.Cells.Item(35, 3).Value = 472.15
